# fix(publipostage): Try to solve Excel emoji problem
# Replace the "statut" emoji values in column A with new text values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of old emoji -> new value (same ordinal position in shared strings):
#   📗 -> ✅
#   📙 -> +3
#   📕 -> -3
#   📘 -> ⚠️

$used = $ws.UsedRange
$rows = $used.Rows.Count

# Column A ("statut") is the only column holding these emoji values.
$col = 1

for ($r = 1; $r -le $rows; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Value2
    if ($val -eq $null) {
        continue
    }
    if ($val -eq "📗") {
        $cell.Value = "✅"
    } elseif ($val -eq "📙") {
        $cell.NumberFormat = "@"
        $cell.Value = "+3"
    } elseif ($val -eq "📕") {
        $cell.NumberFormat = "@"
        $cell.Value = "-3"
    } elseif ($val -eq "📘") {
        $cell.Value = "⚠️"
    }
}
